$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (old C shifts to E, B stays put)
$ws.Range("C:D").EntireColumn.Insert()

# Row 1 headers: B1's old value ("Jun_13") moves to D1, then B1/C1 get the new headers
$ws.Range("D1").Value = $ws.Range("B1").Value()
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# Data rows 2-27: new C and D columns get the same "UN" rating as column B
$ws.Range("C2:D27").Value = "UN"

# Column widths to match the cosmetic formatting change (8.0 "Excel" width units)
$ws.Range("C:C").ColumnWidth = 7.15
$ws.Range("D:D").ColumnWidth = 7.15
$ws.Range("E:E").ColumnWidth = 7.15

Write-Host "done"
